$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column F, shifting the legend (old F/G) to G/H
$ws.Columns("F:F").Insert()

# The column insert doesn't drag the hyperlink anchor along with it, so
# re-create it on the new H7 cell (was G7) pointing at the same target.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H7"), "https://bulbapedia.bulbagarden.net/wiki/Scarlet_%26_Violet_TCG_Series_merchandise")
$ws.Range("H7").Style = "Hyperlink"

# Restyle row 9 (SV 151 / special set) from yellow to green to mark it as released
$ws.Range("A9:D9").Interior.Color = $ws.Range("A2").Interior.Color

# Move the active selection to E12 as in the saved file
$ws.Range("E12").Select()
